# Add "Dataset 9" as a new data column (column J) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new dataset column.
$ws.Cells.Item(1, 10).Value = "Dataset 9"

# Raw values for Dataset 9 (rows 2-26, matching the other dataset columns).
$datasetNine = @{
    2  = 148
    3  = 135
    4  = 233
    5  = 245
    6  = 13
    7  = 170
    8  = 104
    9  = 98
    10 = 277
    11 = 188
    12 = 269
    13 = 201
    14 = 77
    15 = 158
    16 = 269
    17 = 196
    18 = 51
    19 = 79
    20 = 105
    21 = 282
    22 = 16
    23 = 238
    24 = 101
    25 = 265
    26 = 292
}

foreach ($row in $datasetNine.Keys) {
    $ws.Cells.Item($row, 10).Value = $datasetNine[$row]
}

# Leave the cursor parked just past the new data, mirroring the saved selection.
$ws.Range("K1").Select()
